# Apply updated crypto price/volume figures to the "cryptos" worksheet.
# Column D = Price (text-formatted, must stay text even when it looks numeric)
# Column E = Volume(1h) (already protected from numeric coercion by its padding/"%")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.156.04"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.637.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.60%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.96%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.545"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.636.96"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.56%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.145"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.21%  "

$ws.Range("E11").Value = "  -0.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.31%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.351"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.16%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.14%  "

$ws.Range("E15").Value = "  +2.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.116.85"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.71%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.042.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.633.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "363.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.23%  "

$ws.Range("E22").Value = "  +3.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.24%  "

$ws.Range("E24").Value = "  -0.89%  "

$ws.Range("E25").Value = "  +2.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.62%  "

$ws.Range("E28").Value = "  +1.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.769.51"
$ws.Range("D29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "559.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.55%  "

$ws.Range("E33").Value = "  +0.70%  "

$ws.Range("E34").Value = "  +1.08%  "

$ws.Range("E35").Value = "  +1.87%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.04%  "

$ws.Range("E37").Value = "  +3.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.18%  "

$ws.Range("E40").Value = "  +1.44%  "

$ws.Range("E41").Value = "  -0.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₆0342"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.76"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.69%  "

$ws.Range("E45").Value = "  -0.74%  "

$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "158.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.29%  "

$ws.Range("E49").Value = "  +1.77%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.40%  "

$ws.Range("E51").Value = "  +1.31%  "
